# Automatic update of files.
# - Bumps the "Förändrad" (changed) date in column C for every data row
#   from 2023-10-09 (45208) to 2023-10-13 (45212).
# - Rewrites the file-link HYPERLINK() formulas in columns S:Y for the
#   first few rows that have species-find attachments, adding the
#   descriptive filename suffixes (and fixing the "tillsynsmail" folder
#   typo to "ti,llsynsmail") that the upstream generator now emits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Column C: bump every row's changed-date value -----------------
$lastRow = $ws.Cells.Item(1, 1).SpecialCells(11).Row()   # xlCellTypeLastCell = 11

for ($r = 2; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value() -ne $null) {
        $cCell.Value = 45212
    }
}

# --- 2) Columns S:Y: update the hyperlink formulas for rows 2-5 -------
$linkUpdates = @(
    @{ Ref = "S2"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/artfynd/A 30840-2023 artfynd.xlsx'; Label = 'A 30840-2023' },
    @{ Ref = "T2"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/kartor/A 30840-2023 karta.png'; Label = 'A 30840-2023' },
    @{ Ref = "U2"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/knärot/A 30840-2023 karta knärot.png'; Label = 'A 30840-2023' },
    @{ Ref = "V2"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/klagomål/A 30840-2023 fsc-klagomål.docx'; Label = 'A 30840-2023' },
    @{ Ref = "W2"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/klagomålsmail/A 30840-2023 fsc-klagomål mail.docx'; Label = 'A 30840-2023' },
    @{ Ref = "X2"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/tillsyn/A 30840-2023 tillsynsbegäran.docx'; Label = 'A 30840-2023' },
    @{ Ref = "Y2"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/ti,llsynsmail/A 30840-2023 tillsynsbegäran mail.docx'; Label = 'A 30840-2023' },

    @{ Ref = "S3"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/artfynd/A 30841-2023 artfynd.xlsx'; Label = 'A 30841-2023' },
    @{ Ref = "T3"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/kartor/A 30841-2023 karta.png'; Label = 'A 30841-2023' },
    @{ Ref = "U3"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/knärot/A 30841-2023 karta knärot.png'; Label = 'A 30841-2023' },
    @{ Ref = "V3"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/klagomål/A 30841-2023 fsc-klagomål.docx'; Label = 'A 30841-2023' },
    @{ Ref = "W3"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/klagomålsmail/A 30841-2023 fsc-klagomål mail.docx'; Label = 'A 30841-2023' },
    @{ Ref = "X3"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/tillsyn/A 30841-2023 tillsynsbegäran.docx'; Label = 'A 30841-2023' },
    @{ Ref = "Y3"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/ti,llsynsmail/A 30841-2023 tillsynsbegäran mail.docx'; Label = 'A 30841-2023' },

    @{ Ref = "S4"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/artfynd/A 30839-2023 artfynd.xlsx'; Label = 'A 30839-2023' },
    @{ Ref = "T4"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/kartor/A 30839-2023 karta.png'; Label = 'A 30839-2023' },
    @{ Ref = "V4"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/klagomål/A 30839-2023 fsc-klagomål.docx'; Label = 'A 30839-2023' },
    @{ Ref = "W4"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/klagomålsmail/A 30839-2023 fsc-klagomål mail.docx'; Label = 'A 30839-2023' },
    @{ Ref = "X4"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/tillsyn/A 30839-2023 tillsynsbegäran.docx'; Label = 'A 30839-2023' },
    @{ Ref = "Y4"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2260/ti,llsynsmail/A 30839-2023 tillsynsbegäran mail.docx'; Label = 'A 30839-2023' },

    @{ Ref = "S5"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2283/artfynd/A 33036-2023 artfynd.xlsx'; Label = 'A 33036-2023' },
    @{ Ref = "T5"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2283/kartor/A 33036-2023 karta.png'; Label = 'A 33036-2023' },
    @{ Ref = "V5"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2283/klagomål/A 33036-2023 fsc-klagomål.docx'; Label = 'A 33036-2023' },
    @{ Ref = "W5"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2283/klagomålsmail/A 33036-2023 fsc-klagomål mail.docx'; Label = 'A 33036-2023' },
    @{ Ref = "X5"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2283/tillsyn/A 33036-2023 tillsynsbegäran.docx'; Label = 'A 33036-2023' },
    @{ Ref = "Y5"; Url = 'https://klasma.github.io/LoggingDetectiveFiles/Logging_2283/ti,llsynsmail/A 33036-2023 tillsynsbegäran mail.docx'; Label = 'A 33036-2023' }
)

foreach ($upd in $linkUpdates) {
    $cell = $ws.Range($upd.Ref)
    $cell.Formula = '=HYPERLINK("' + $upd.Url + '", "' + $upd.Label + '")'
}
